$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet used to start with a big "header" row (row 1) and a blank
# spacer row (row 2) before the real tabular data in row 3. Those two
# rows were being stripped by mistake - undo that by deleting them
# outright so the real data (old row 3) becomes row 1, with everything
# below shifting up to match.
$ws.Rows("1:2").Delete()

# Reset the active selection back to the top-left cell.
$ws.Range("A1").Select()
